# Commit: "Added data from 3/26/2020"
# Appends 91 new case rows (rows 366-456) to the houstonNumbers sheet, one
# per new COVID-19 case reported on 3/26/2020 (Excel date serial 43916),
# with only the Date (col A) and County (col D) populated - matching the
# existing sparse data pattern used for newly-added days in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 43916

# (first row, last row, county) for each contiguous run of same-county cases
$countyRuns = @(
    @{start=366; end=371; county="Brazoria"},
    @{start=372; end=403; county="Fort Bend"},
    @{start=404; end=416; county="Galveston"},
    @{start=417; end=418; county="Washington"},
    @{start=419; end=428; county="Montgomery"},
    @{start=429; end=444; county="Harris"},
    @{start=445; end=446; county="Houston"},
    @{start=447; end=449; county="Walker"},
    @{start=450; end=456; county="Brazos"}
)

foreach ($run in $countyRuns) {
    for ($r = $run.start; $r -le $run.end; $r++) {
        # Copy the date formatting (short-date number format) from an
        # existing date cell so the new cells share the same style index
        # instead of Excel minting a brand-new one.
        $ws.Range("A2").Copy()
        $ws.Range("A$r").PasteSpecial(-4122)

        $ws.Range("A$r").Value = $newDate
        $ws.Range("D$r").Value = $run.county
    }
}

# Leave the selection where the author's session ended up.
[void]$ws.Range("C439").Select()
$excel.CutCopyMode = $false
